$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4061.25
$ws.Range("I62").Value = 4250
$ws.Range("J62").Value = 3998.3333
$ws.Range("K62").Value = 4250
$ws.Range("L62").Value = 3998.3333
$ws.Range("M62").Value = -3626
$ws.Range("N62").Value = -5246.3333
$ws.Range("H65").Value = 4061.25
$ws.Range("I65").Value = 4250
$ws.Range("J65").Value = 3998.3333
$ws.Range("K65").Value = 21250
$ws.Range("L65").Value = 19991.6665
$ws.Range("M65").Value = -18130
$ws.Range("N65").Value = -26231.6665
$ws.Range("H100").Value = 3331.1667
$ws.Range("I100").Value = 2499.25
$ws.Range("K100").Value = 2499.25
$ws.Range("M100").Value = -1958.25
$ws.Range("H137").Value = 1172
$ws.Range("I137").Value = 1021.0769
$ws.Range("J137").Value = 1662.5
$ws.Range("K137").Value = 3063.2307
$ws.Range("L137").Value = 4987.5
$ws.Range("M137").Value = -513.2307000000001
$ws.Range("N137").Value = -10087.5
$ws.Range("H138").Value = 4258.8906
$ws.Range("I138").Value = 1329.2084
$ws.Range("J138").Value = 5693.837
$ws.Range("K138").Value = 3987.6252
$ws.Range("L138").Value = 17081.511
$ws.Range("M138").Value = 1152.3748
$ws.Range("N138").Value = -27361.511

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 43596430
$ws.Range("I132").Value = 20997.357
$ws.Range("J132").Value = 111380440
$ws.Range("K132").Value = 62992.071
$ws.Range("L132").Value = 334141320
$ws.Range("M132").Value = -60462.071
$ws.Range("N132").Value = -334146380

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0
$ws.Range("H96").Value = 36373.875
$ws.Range("I96").Value = 8713
$ws.Range("J96").Value = 230000
$ws.Range("K96").Value = 8713
$ws.Range("L96").Value = 230000
$ws.Range("M96").Value = -5967
$ws.Range("N96").Value = -235492
$ws.Range("H99").Value = 3702.9285
$ws.Range("I99").Value = 3334.3
$ws.Range("K99").Value = 3334.3
$ws.Range("M99").Value = -1836.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0
$ws.Range("H132").Value = 103887.3
$ws.Range("I132").Value = 253570.62
$ws.Range("J132").Value = 4098.4165
$ws.Range("K132").Value = 760711.86
$ws.Range("L132").Value = 12295.2495
$ws.Range("M132").Value = -758181.86
$ws.Range("N132").Value = -17355.2495
$ws.Range("H135").Value = 199998
$ws.Range("J135").Value = 199998
$ws.Range("L135").Value = 199998
$ws.Range("N135").Value = -210138

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34312656
$ws.Range("I4").Value = 41715176
$ws.Range("J4").Value = 24182894
$ws.Range("K4").Value = 125145528
$ws.Range("L4").Value = 72548682
$ws.Range("M4").Value = -125145416
$ws.Range("N4").Value = -72548906
$ws.Range("H68").Value = 948
$ws.Range("I68").Value = 997.6667
$ws.Range("J68").Value = 799
$ws.Range("K68").Value = 2993.0001
$ws.Range("L68").Value = 2397
$ws.Range("M68").Value = -2182.0001
$ws.Range("N68").Value = -4019
$ws.Range("H71").Value = 948
$ws.Range("I71").Value = 997.6667
$ws.Range("J71").Value = 799
$ws.Range("K71").Value = 8979.0003
$ws.Range("L71").Value = 7191
$ws.Range("M71").Value = -4923.0003
$ws.Range("N71").Value = -15303
$ws.Range("H133").Value = 1806.4286
$ws.Range("J133").Value = 3033
$ws.Range("L133").Value = 9099
$ws.Range("N133").Value = -19219
$ws.Range("H134").Value = 1274.5714
$ws.Range("I134").Value = 981.5
$ws.Range("J134").Value = 3033
$ws.Range("K134").Value = 2944.5
$ws.Range("L134").Value = 9099
$ws.Range("M134").Value = 2125.5
$ws.Range("N134").Value = -19239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 160.66667
$ws.Range("I46").Value = 160.66667
$ws.Range("K46").Value = 160.66667
$ws.Range("M46").Value = -4.666670000000011
$ws.Range("H80").Value = 5981.5137
$ws.Range("I80").Value = 5366.864
$ws.Range("J80").Value = 6883
$ws.Range("K80").Value = 5366.864
$ws.Range("L80").Value = 6883
$ws.Range("M80").Value = -4368.864
$ws.Range("N80").Value = -8879
$ws.Range("H83").Value = 5981.5137
$ws.Range("I83").Value = 5366.864
$ws.Range("J83").Value = 6883
$ws.Range("K83").Value = 26834.32
$ws.Range("L83").Value = 34415
$ws.Range("M83").Value = -21842.32
$ws.Range("N83").Value = -44399
$ws.Range("H97").Value = 862
$ws.Range("I97").Value = 880.8095
$ws.Range("J97").Value = 796.1667
$ws.Range("K97").Value = 880.8095
$ws.Range("L97").Value = 796.1667
$ws.Range("M97").Value = -384.8095
$ws.Range("N97").Value = -1788.1667
$ws.Range("H113").Value = 4209
$ws.Range("I113").Value = 2563.75
$ws.Range("K113").Value = 2563.75
$ws.Range("M113").Value = -393.75
$ws.Range("H135").Value = 183322.5
$ws.Range("J135").Value = 183322.5
$ws.Range("L135").Value = 183322.5
$ws.Range("N135").Value = -193462.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3180.5625
$ws.Range("I7").Value = 2407.4167
$ws.Range("J7").Value = 5500
$ws.Range("K7").Value = 2407.4167
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = -2295.4167
$ws.Range("N7").Value = -5724
$ws.Range("H25").Value = 2401.6
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 3336
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 3336
$ws.Range("M25").Value = -770
$ws.Range("N25").Value = -3796
$ws.Range("H61").Value = 4416.364
$ws.Range("I61").Value = 3450.125
$ws.Range("J61").Value = 6993
$ws.Range("K61").Value = 3450.125
$ws.Range("L61").Value = 6993
$ws.Range("M61").Value = -3248.125
$ws.Range("N61").Value = -7397
$ws.Range("H68").Value = 2519.111
$ws.Range("I68").Value = 2467.5715
$ws.Range("K68").Value = 2467.5715
$ws.Range("M68").Value = -1718.5715
$ws.Range("H71").Value = 2519.111
$ws.Range("I71").Value = 2467.5715
$ws.Range("K71").Value = 12337.8575
$ws.Range("M71").Value = -8593.8575
$ws.Range("H100").Value = 3696.8
$ws.Range("I100").Value = 3190.889
$ws.Range("K100").Value = 3190.889
$ws.Range("M100").Value = -2649.889
$ws.Range("H113").Value = 4416.364
$ws.Range("I113").Value = 3450.125
$ws.Range("J113").Value = 6993
$ws.Range("K113").Value = 3450.125
$ws.Range("L113").Value = 6993
$ws.Range("M113").Value = -1280.125
$ws.Range("N113").Value = -11333
$ws.Range("H126").Value = 3180.5625
$ws.Range("I126").Value = 2407.4167
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 7222.250100000001
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -4752.250100000001
$ws.Range("N126").Value = -21440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1516.25
$ws.Range("I107").Value = 613.5
$ws.Range("K107").Value = 1840.5
$ws.Range("M107").Value = 79.5
$ws.Range("H141").Value = 68219.89
$ws.Range("J141").Value = 63122.375
$ws.Range("L141").Value = 63122.375
$ws.Range("N141").Value = -73482.375
